$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("a Scarlett devkit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "an Xbox Series X|S devkit"
